$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.424.09"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("D3").Value = "1.860.95"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'325.26"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("D6").Value = "'1.006"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "'0.4553"
$ws.Range("E7").Value = "  -2.28%  "
$ws.Range("D8").Value = "'0.3827"
$ws.Range("E8").Value = "  -2.24%  "
$ws.Range("D9").Value = "'0.07803"
$ws.Range("E9").Value = "  -1.19%  "
$ws.Range("D10").Value = "'0.9842"
$ws.Range("E10").Value = "  +1.05%  "
$ws.Range("D11").Value = "'21.45"
$ws.Range("E11").Value = "  -3.67%  "
$ws.Range("D12").Value = "1.843.77"
$ws.Range("E12").Value = "  -2.58%  "
$ws.Range("D13").Value = "'5.630"
$ws.Range("E13").Value = "  -1.83%  "
$ws.Range("D14").Value = "'6.883"
$ws.Range("E14").Value = "  -0.99%  "
$ws.Range("D15").Value = "'0.06914"
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").Value = "'86.39"
$ws.Range("E17").Value = "  -2.77%  "
$ws.Range("D18").Value = "'0.000009928"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("D19").Value = "'16.64"
$ws.Range("E19").Value = "  -1.62%  "
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("D21").Value = "28.437.89"
$ws.Range("E21").Value = "  -0.75%  "
$ws.Range("D22").Value = "'5.244"
$ws.Range("E22").Value = "  -1.60%  "
$ws.Range("D23").Value = "'10.86"
$ws.Range("E23").Value = "  -1.86%  "
$ws.Range("D24").Value = "'2.090"
$ws.Range("E24").Value = "  -1.76%  "
$ws.Range("D25").Value = "2.068.91"
$ws.Range("E25").Value = "  -1.99%  "
$ws.Range("D26").Value = "'153.12"
$ws.Range("E26").Value = "  -1.31%  "
$ws.Range("D27").Value = "'19.09"
$ws.Range("E27").Value = "  -1.05%  "
$ws.Range("D28").Value = "'5.642"
$ws.Range("E28").Value = "  -2.59%  "
$ws.Range("D29").Value = "'117.16"
$ws.Range("E29").Value = "  -1.78%  "
$ws.Range("D30").Value = "'1.895"
$ws.Range("E30").Value = "  -4.93%  "
$ws.Range("D31").Value = "'0.09276"
$ws.Range("E31").Value = "  -0.89%  "
$ws.Range("D32").Value = "'0.9043"
$ws.Range("E32").Value = "  -3.71%  "
$ws.Range("D33").Value = "'5.254"
$ws.Range("D34").Value = "'1.312"
$ws.Range("E34").Value = "  -2.22%  "
$ws.Range("D35").Value = "'3.289"
$ws.Range("E35").Value = "  -1.68%  "
$ws.Range("D36").Value = "'0.05678"
$ws.Range("E36").Value = "  -2.80%  "
$ws.Range("D37").Value = "'1.151"
$ws.Range("E37").Value = "  -0.31%  "
$ws.Range("D38").Value = "'0.02029"
$ws.Range("E38").Value = "  -4.05%  "
$ws.Range("D39").Value = "'7.633"
$ws.Range("E39").Value = "  -3.29%  "
$ws.Range("D40").Value = "'0.5537"
$ws.Range("E40").Value = "  -2.00%  "
$ws.Range("D41").Value = "'0.1761"
$ws.Range("E41").Value = "  -0.87%  "
$ws.Range("D42").Value = "'9.590"
$ws.Range("E42").Value = "  -3.83%  "
$ws.Range("D43").Value = "'0.07117"
$ws.Range("E43").Value = "  -3.36%  "
$ws.Range("D44").Value = "'11.55"
$ws.Range("E44").Value = "  -1.09%  "
$ws.Range("D45").Value = "'0.5219"
$ws.Range("E45").Value = "  -2.01%  "
$ws.Range("D46").Value = "'1.124"
$ws.Range("E46").Value = "  -1.41%  "
$ws.Range("D47").Value = "'2.096"
$ws.Range("E47").Value = "  -4.38%  "
$ws.Range("D48").Value = "'1.802"
$ws.Range("E48").Value = "  -2.65%  "
$ws.Range("D49").Value = "'111.64"
$ws.Range("E49").Value = "  -2.18%  "
$ws.Range("D50").Value = "'2.436"
$ws.Range("E50").Value = "  +3.43%  "
$ws.Range("D51").Value = "'1.006"
$ws.Range("E51").Value = "  +0.10%  "
